$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting the old E:I block to F:J
$ws.Range("E1").EntireColumn.Insert()

# Give the new trailing header cell (K1) the same formatting as its
# neighbour (J1, bold/bordered/centered header style) before filling it in.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Header row
$ws.Range("E1").Value = "1X2"
$ws.Range("G1").Value = "RTDO L.1"
$ws.Range("H1").Value = "LOCAL.1"
$ws.Range("I1").Value = "VISITANTE.1"
$ws.Range("J1").Value = "RTDO V.1"
$ws.Range("K1").Value = "1X2.1"

# Row 2
$ws.Range("A2").Value = 43
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 36
$ws.Range("H2").Value = "Puche"
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 1

# Row 3
$ws.Range("A3").Value = 50
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 31
$ws.Range("J3").Value = 49
$ws.Range("K3").Value = 2

# Row 4
$ws.Range("A4").Value = 41
$ws.Range("D4").Value = 23
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = 38
$ws.Range("J4").Value = 32
$ws.Range("K4").Value = 1
